# Updated cryptos list on Wed Oct 18 05:08:03 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.703.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.65%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.17%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.492"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.61"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.80%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.22%  "

# Row 10
$ws.Range("E10").Value = "  -1.55%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0592"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0889"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.796.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.570.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.96%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.522"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.78%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.674.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.40%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.67%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.55%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.33%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0693"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.32%  "

# Row 22
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.95%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.89%  "

# Row 25
$ws.Range("E25").Value = "  +10.61%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.16%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.94%  "

# Row 29
$ws.Range("E29").Value = "  -3.01%  "

# Row 30
$ws.Range("E30").Value = "  -0.11%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0487"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.78%  "

# Row 32
$ws.Range("E32").Value = "  -2.40%  "

# Row 33
$ws.Range("E33").Value = "  -0.71%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.52%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.393.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.70%  "

# Row 36
$ws.Range("E36").Value = "  +1.53%  "

# Row 37
$ws.Range("E37").Value = "  -3.16%  "

# Row 38
$ws.Range("E38").Value = "  +0.62%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.80%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0167"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.528"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.66%  "

# Row 42
$ws.Range("E42").Value = "  -0.07%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.35%  "

# Row 44
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.793"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.66%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0469"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.05%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.966"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.70%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.710.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.25%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.96%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0518"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.27%  "
